$wb = $excel.ActiveWorkbook

# --- Sheet 1: compare_models ---
$ws = $wb.Worksheets.Item("compare_models")
$ws.Range("A2").Value = "et"
$ws.Range("B2").Value = "Extra Trees Regressor"
$ws.Range("C2").Value = 5.1426
$ws.Range("D2").Value = 56.5649
$ws.Range("E2").Value = 7.2851
$ws.Range("F2").Value = 0.9261
$ws.Range("G2").Value = 0.1268
$ws.Range("H2").Value = 0.09379999999999999
$ws.Range("I2").Value = 0.106

$ws.Range("A3").Value = "lightgbm"
$ws.Range("B3").Value = "Light Gradient Boosting Machine"
$ws.Range("C3").Value = 6.4696
$ws.Range("D3").Value = 78.3138
$ws.Range("E3").Value = 8.7201
$ws.Range("F3").Value = 0.8969
$ws.Range("G3").Value = 0.153
$ws.Range("H3").Value = 0.1197
$ws.Range("I3").Value = 0.036

$ws.Range("A4").Value = "gbr"
$ws.Range("B4").Value = "Gradient Boosting Regressor"
$ws.Range("C4").Value = 6.2325
$ws.Range("D4").Value = 84.2282
$ws.Range("E4").Value = 8.9201
$ws.Range("F4").Value = 0.8908
$ws.Range("G4").Value = 0.1492
$ws.Range("H4").Value = 0.1113
$ws.Range("I4").Value = 0.062

$ws.Range("A5").Value = "rf"
$ws.Range("B5").Value = "Random Forest Regressor"
$ws.Range("C5").Value = 6.3502
$ws.Range("D5").Value = 88.59869999999999
$ws.Range("E5").Value = 9.1892
$ws.Range("F5").Value = 0.8847
$ws.Range("G5").Value = 0.1509
$ws.Range("H5").Value = 0.113
$ws.Range("I5").Value = 0.14

$ws.Range("A6").Value = "ada"
$ws.Range("B6").Value = "AdaBoost Regressor"
$ws.Range("C6").Value = 7.0161
$ws.Range("D6").Value = 92.4585
$ws.Range("E6").Value = 9.286199999999999
$ws.Range("F6").Value = 0.8794
$ws.Range("G6").Value = 0.1588
$ws.Range("H6").Value = 0.1293
$ws.Range("I6").Value = 0.062

$ws.Range("A7").Value = "lr"
$ws.Range("B7").Value = "Linear Regression"
$ws.Range("C7").Value = 8.07
$ws.Range("D7").Value = 108.3445
$ws.Range("E7").Value = 10.2559
$ws.Range("F7").Value = 0.8574000000000001
$ws.Range("G7").Value = 0.2263
$ws.Range("H7").Value = 0.1565
$ws.Range("I7").Value = 1.894

$ws.Range("A8").Value = "knn"
$ws.Range("B8").Value = "K Neighbors Regressor"
$ws.Range("C8").Value = 7.0335
$ws.Range("D8").Value = 111.519
$ws.Range("E8").Value = 10.3501
$ws.Range("F8").Value = 0.8515
$ws.Range("G8").Value = 0.1723
$ws.Range("H8").Value = 0.1302
$ws.Range("I8").Value = 0.03

$ws.Range("A9").Value = "ridge"
$ws.Range("B9").Value = "Ridge Regression"
$ws.Range("C9").Value = 8.3019
$ws.Range("D9").Value = 114.0235
$ws.Range("E9").Value = 10.5637
$ws.Range("F9").Value = 0.8498
$ws.Range("G9").Value = 0.2042
$ws.Range("H9").Value = 0.1594
$ws.Range("I9").Value = 0.022

$ws.Range("A10").Value = "br"
$ws.Range("B10").Value = "Bayesian Ridge"
$ws.Range("C10").Value = 8.312099999999999
$ws.Range("D10").Value = 117.8118
$ws.Range("E10").Value = 10.7621
$ws.Range("F10").Value = 0.8443000000000001
$ws.Range("G10").Value = 0.1965
$ws.Range("H10").Value = 0.1584
$ws.Range("I10").Value = 0.022

$ws.Range("A11").Value = "lasso"
$ws.Range("B11").Value = "Lasso Regression"
$ws.Range("C11").Value = 8.315899999999999
$ws.Range("D11").Value = 133.3326
$ws.Range("E11").Value = 11.4342
$ws.Range("F11").Value = 0.8235
$ws.Range("G11").Value = 0.1939
$ws.Range("H11").Value = 0.1558
$ws.Range("I11").Value = 0.748

$ws.Range("A12").Value = "en"
$ws.Range("B12").Value = "Elastic Net"
$ws.Range("C12").Value = 8.651999999999999
$ws.Range("D12").Value = 144.5511
$ws.Range("E12").Value = 11.9127
$ws.Range("F12").Value = 0.8084
$ws.Range("G12").Value = 0.2025
$ws.Range("H12").Value = 0.1633
$ws.Range("I12").Value = 0.02

$ws.Range("A13").Value = "huber"
$ws.Range("B13").Value = "Huber Regressor"
$ws.Range("C13").Value = 8.994300000000001
$ws.Range("D13").Value = 146.9392
$ws.Range("E13").Value = 11.9851
$ws.Range("F13").Value = 0.8023
$ws.Range("G13").Value = 0.2589
$ws.Range("H13").Value = 0.1741
$ws.Range("I13").Value = 0.034

$ws.Range("A14").Value = "par"
$ws.Range("B14").Value = "Passive Aggressive Regressor"
$ws.Range("C14").Value = 9.949400000000001
$ws.Range("D14").Value = 171.2201
$ws.Range("E14").Value = 12.9285
$ws.Range("F14").Value = 0.7771
$ws.Range("G14").Value = 0.2649
$ws.Range("H14").Value = 0.1864
$ws.Range("I14").Value = 0.028

$ws.Range("A15").Value = "dt"
$ws.Range("B15").Value = "Decision Tree Regressor"
$ws.Range("C15").Value = 9.1456
$ws.Range("D15").Value = 191.6835
$ws.Range("E15").Value = 13.4555
$ws.Range("F15").Value = 0.7514999999999999
$ws.Range("G15").Value = 0.2188
$ws.Range("H15").Value = 0.1592
$ws.Range("I15").Value = 0.03

$ws.Range("A16").Value = "omp"
$ws.Range("B16").Value = "Orthogonal Matching Pursuit"
$ws.Range("C16").Value = 9.303800000000001
$ws.Range("D16").Value = 187.351
$ws.Range("E16").Value = 13.4743
$ws.Range("F16").Value = 0.751
$ws.Range("G16").Value = 0.2338
$ws.Range("H16").Value = 0.1769
$ws.Range("I16").Value = 0.028

$ws.Range("A17").Value = "llar"
$ws.Range("B17").Value = "Lasso Least Angle Regression"
$ws.Range("C17").Value = 14.8057
$ws.Range("D17").Value = 328.4843
$ws.Range("E17").Value = 17.9968
$ws.Range("F17").Value = 0.5672
$ws.Range("G17").Value = 0.3218
$ws.Range("H17").Value = 0.3076
$ws.Range("I17").Value = 0.024

$ws.Range("A18").Value = "dummy"
$ws.Range("B18").Value = "Dummy Regressor"
$ws.Range("C18").Value = 23.6246
$ws.Range("D18").Value = 760.3468
$ws.Range("E18").Value = 27.512
$ws.Range("F18").Value = -0.0098
$ws.Range("G18").Value = 0.483
$ws.Range("H18").Value = 0.5036
$ws.Range("I18").Value = 0.032

$ws.Range("A19").Value = "lar"
$ws.Range("B19").Value = "Least Angle Regression"
$ws.Range("C19").Value = 28.8631
$ws.Range("D19").Value = 1675.6423
$ws.Range("E19").Value = 36.5242
$ws.Range("F19").Value = -1.0545
$ws.Range("G19").Value = 0.6599
$ws.Range("H19").Value = 0.61
$ws.Range("I19").Value = 0.022

# --- Sheet: tuned_1 ---
$ws = $wb.Worksheets.Item("tuned_1")
$ws.Range("B2").Value = 4.7341
$ws.Range("C2").Value = 54.2063
$ws.Range("D2").Value = 7.3625
$ws.Range("E2").Value = 0.9249000000000001
$ws.Range("F2").Value = 0.1171
$ws.Range("G2").Value = 0.07920000000000001

$ws.Range("B3").Value = 5.5477
$ws.Range("C3").Value = 68.6525
$ws.Range("D3").Value = 8.2857
$ws.Range("E3").Value = 0.9184
$ws.Range("F3").Value = 0.1292
$ws.Range("G3").Value = 0.094

$ws.Range("B4").Value = 4.474
$ws.Range("C4").Value = 34.2661
$ws.Range("D4").Value = 5.8537
$ws.Range("E4").Value = 0.9412
$ws.Range("F4").Value = 0.1204
$ws.Range("G4").Value = 0.0916

$ws.Range("B5").Value = 5.1587
$ws.Range("C5").Value = 50.2205
$ws.Range("D5").Value = 7.0866
$ws.Range("E5").Value = 0.9405
$ws.Range("F5").Value = 0.111
$ws.Range("G5").Value = 0.08359999999999999

$ws.Range("B6").Value = 7.6564
$ws.Range("C6").Value = 121.156
$ws.Range("D6").Value = 11.0071
$ws.Range("E6").Value = 0.8433
$ws.Range("F6").Value = 0.2116
$ws.Range("G6").Value = 0.1644

$ws.Range("B7").Value = 5.5142
$ws.Range("C7").Value = 65.7003
$ws.Range("D7").Value = 7.9191
$ws.Range("E7").Value = 0.9137
$ws.Range("F7").Value = 0.1378
$ws.Range("G7").Value = 0.1026

$ws.Range("B8").Value = 1.132
$ws.Range("C8").Value = 29.8124
$ws.Range("D8").Value = 1.7285
$ws.Range("E8").Value = 0.0363
$ws.Range("F8").Value = 0.0373
$ws.Range("G8").Value = 0.0314

# --- Sheet: tuned_2 ---
$ws = $wb.Worksheets.Item("tuned_2")
$ws.Range("B2").Value = 6.1607
$ws.Range("C2").Value = 64.52760000000001
$ws.Range("D2").Value = 8.0329
$ws.Range("E2").Value = 0.9106
$ws.Range("F2").Value = 0.1292
$ws.Range("G2").Value = 0.1073

$ws.Range("B3").Value = 5.9446
$ws.Range("C3").Value = 73.49420000000001
$ws.Range("D3").Value = 8.572900000000001
$ws.Range("E3").Value = 0.9126
$ws.Range("F3").Value = 0.1373
$ws.Range("G3").Value = 0.1042

$ws.Range("B4").Value = 4.3214
$ws.Range("C4").Value = 34.1426
$ws.Range("D4").Value = 5.8432
$ws.Range("E4").Value = 0.9414
$ws.Range("F4").Value = 0.1045
$ws.Range("G4").Value = 0.08160000000000001

$ws.Range("B5").Value = 5.7123
$ws.Range("C5").Value = 48.128
$ws.Range("D5").Value = 6.9374
$ws.Range("E5").Value = 0.9429999999999999
$ws.Range("F5").Value = 0.1087
$ws.Range("G5").Value = 0.09180000000000001

$ws.Range("B6").Value = 7.0208
$ws.Range("C6").Value = 92.6619
$ws.Range("D6").Value = 9.626099999999999
$ws.Range("E6").Value = 0.8802
$ws.Range("F6").Value = 0.1935
$ws.Range("G6").Value = 0.1505

$ws.Range("B7").Value = 5.832
$ws.Range("C7").Value = 62.5909
$ws.Range("D7").Value = 7.8025
$ws.Range("E7").Value = 0.9176
$ws.Range("F7").Value = 0.1346
$ws.Range("G7").Value = 0.1071

$ws.Range("B8").Value = 0.8753
$ws.Range("C8").Value = 20.2256
$ws.Range("D8").Value = 1.3084
$ws.Range("E8").Value = 0.0232
$ws.Range("F8").Value = 0.0319
$ws.Range("G8").Value = 0.0236

# --- Sheet: tuned_3 ---
$ws = $wb.Worksheets.Item("tuned_3")
$ws.Range("B2").Value = 4.8012
$ws.Range("C2").Value = 43.7897
$ws.Range("D2").Value = 6.6174
$ws.Range("E2").Value = 0.9393
$ws.Range("F2").Value = 0.09909999999999999
$ws.Range("G2").Value = 0.0798

$ws.Range("B3").Value = 5.0161
$ws.Range("C3").Value = 54.1704
$ws.Range("D3").Value = 7.3601
$ws.Range("E3").Value = 0.9356
$ws.Range("F3").Value = 0.1212
$ws.Range("G3").Value = 0.0915

$ws.Range("B4").Value = 3.8946
$ws.Range("C4").Value = 27.9737
$ws.Range("D4").Value = 5.289
$ws.Range("E4").Value = 0.952
$ws.Range("F4").Value = 0.0901
$ws.Range("G4").Value = 0.073

$ws.Range("B5").Value = 5.0757
$ws.Range("C5").Value = 46.3666
$ws.Range("D5").Value = 6.8093
$ws.Range("E5").Value = 0.9451000000000001
$ws.Range("F5").Value = 0.1223
$ws.Range("G5").Value = 0.0843

$ws.Range("B6").Value = 7.7388
$ws.Range("C6").Value = 104.4714
$ws.Range("D6").Value = 10.2211
$ws.Range("E6").Value = 0.8649
$ws.Range("F6").Value = 0.2074
$ws.Range("G6").Value = 0.1682

$ws.Range("B7").Value = 5.3053
$ws.Range("C7").Value = 55.3544
$ws.Range("D7").Value = 7.2594
$ws.Range("E7").Value = 0.9274
$ws.Range("F7").Value = 0.128
$ws.Range("G7").Value = 0.0994

$ws.Range("B8").Value = 1.2886
$ws.Range("C8").Value = 25.9927
$ws.Range("D8").Value = 1.6297
$ws.Range("E8").Value = 0.0317
$ws.Range("F8").Value = 0.0416
$ws.Range("G8").Value = 0.0349

# --- Sheet: tuned_4 ---
$ws = $wb.Worksheets.Item("tuned_4")
$ws.Range("B2").Value = 5.1418
$ws.Range("C2").Value = 64.4306
$ws.Range("D2").Value = 8.026899999999999
$ws.Range("E2").Value = 0.9107
$ws.Range("F2").Value = 0.1243
$ws.Range("G2").Value = 0.0852

$ws.Range("B3").Value = 5.7532
$ws.Range("C3").Value = 72.5934
$ws.Range("D3").Value = 8.520200000000001
$ws.Range("E3").Value = 0.9137
$ws.Range("F3").Value = 0.1259
$ws.Range("G3").Value = 0.0927

$ws.Range("B4").Value = 4.607
$ws.Range("C4").Value = 37.2248
$ws.Range("D4").Value = 6.1012
$ws.Range("E4").Value = 0.9361
$ws.Range("F4").Value = 0.1167
$ws.Range("G4").Value = 0.0916

$ws.Range("B5").Value = 5.9799
$ws.Range("C5").Value = 64.8425
$ws.Range("D5").Value = 8.0525
$ws.Range("E5").Value = 0.9232
$ws.Range("F5").Value = 0.1258
$ws.Range("G5").Value = 0.09569999999999999

$ws.Range("B6").Value = 7.779
$ws.Range("C6").Value = 137.6488
$ws.Range("D6").Value = 11.7324
$ws.Range("E6").Value = 0.822
$ws.Range("F6").Value = 0.2167
$ws.Range("G6").Value = 0.1669

$ws.Range("B7").Value = 5.8522
$ws.Range("C7").Value = 75.348
$ws.Range("D7").Value = 8.486599999999999
$ws.Range("E7").Value = 0.9012
$ws.Range("F7").Value = 0.1419
$ws.Range("G7").Value = 0.1064

$ws.Range("B8").Value = 1.0765
$ws.Range("C8").Value = 33.3824
$ws.Range("D8").Value = 1.8235
$ws.Range("E8").Value = 0.0406
$ws.Range("F8").Value = 0.0376
$ws.Range("G8").Value = 0.0304

# --- Sheet: tuned_5 ---
$ws = $wb.Worksheets.Item("tuned_5")
$ws.Range("B2").Value = 6.5565
$ws.Range("C2").Value = 72.0138
$ws.Range("D2").Value = 8.4861
$ws.Range("E2").Value = 0.9002
$ws.Range("F2").Value = 0.1531
$ws.Range("G2").Value = 0.1259

$ws.Range("B3").Value = 6.5134
$ws.Range("C3").Value = 77.18770000000001
$ws.Range("D3").Value = 8.7857
$ws.Range("E3").Value = 0.9083
$ws.Range("F3").Value = 0.137
$ws.Range("G3").Value = 0.1078

$ws.Range("B4").Value = 5.3312
$ws.Range("C4").Value = 47.8971
$ws.Range("D4").Value = 6.9208
$ws.Range("E4").Value = 0.9177999999999999
$ws.Range("F4").Value = 0.1296
$ws.Range("G4").Value = 0.1049

$ws.Range("B5").Value = 6.8055
$ws.Range("C5").Value = 78.1621
$ws.Range("D5").Value = 8.8409
$ws.Range("E5").Value = 0.9074
$ws.Range("F5").Value = 0.147
$ws.Range("G5").Value = 0.113

$ws.Range("B6").Value = 8.855
$ws.Range("C6").Value = 160.6664
$ws.Range("D6").Value = 12.6754
$ws.Range("E6").Value = 0.7922
$ws.Range("F6").Value = 0.227
$ws.Range("G6").Value = 0.1853

$ws.Range("B7").Value = 6.8123
$ws.Range("C7").Value = 87.1854
$ws.Range("D7").Value = 9.1418
$ws.Range("E7").Value = 0.8852
$ws.Range("F7").Value = 0.1587
$ws.Range("G7").Value = 0.1274

$ws.Range("B8").Value = 1.142
$ws.Range("C8").Value = 38.3525
$ws.Range("D8").Value = 1.9009
$ws.Range("E8").Value = 0.0468
$ws.Range("F8").Value = 0.0351
$ws.Range("G8").Value = 0.0298

# --- Sheet: blend_model ---
$ws = $wb.Worksheets.Item("blend_model")
$ws.Range("B2").Value = 4.7905
$ws.Range("C2").Value = 46.6769
$ws.Range("D2").Value = 6.832
$ws.Range("E2").Value = 0.9353
$ws.Range("F2").Value = 0.1089
$ws.Range("G2").Value = 0.08260000000000001

$ws.Range("B3").Value = 5.3063
$ws.Range("C3").Value = 56.9215
$ws.Range("D3").Value = 7.5446
$ws.Range("E3").Value = 0.9323
$ws.Range("F3").Value = 0.1195
$ws.Range("G3").Value = 0.089

$ws.Range("B4").Value = 3.7327
$ws.Range("C4").Value = 25.842
$ws.Range("D4").Value = 5.0835
$ws.Range("E4").Value = 0.9557
$ws.Range("F4").Value = 0.09950000000000001
$ws.Range("G4").Value = 0.0762

$ws.Range("B5").Value = 5.2837
$ws.Range("C5").Value = 49.3211
$ws.Range("D5").Value = 7.0229
$ws.Range("E5").Value = 0.9416
$ws.Range("F5").Value = 0.1145
$ws.Range("G5").Value = 0.0863

$ws.Range("B6").Value = 7.3288
$ws.Range("C6").Value = 112.4465
$ws.Range("D6").Value = 10.6041
$ws.Range("E6").Value = 0.8546
$ws.Range("F6").Value = 0.2025
$ws.Range("G6").Value = 0.1561

$ws.Range("B7").Value = 5.2884
$ws.Range("C7").Value = 58.2416
$ws.Range("D7").Value = 7.4174
$ws.Range("E7").Value = 0.9239000000000001
$ws.Range("F7").Value = 0.129
$ws.Range("G7").Value = 0.098

$ws.Range("B8").Value = 1.1689
$ws.Range("C8").Value = 28.993
$ws.Range("D8").Value = 1.7954
$ws.Range("E8").Value = 0.0356
$ws.Range("F8").Value = 0.0373
$ws.Range("G8").Value = 0.0293

# --- Sheet: stack_model ---
$ws = $wb.Worksheets.Item("stack_model")
$ws.Range("B2").Value = 4.8726
$ws.Range("C2").Value = 45.4692
$ws.Range("D2").Value = 6.7431
$ws.Range("E2").Value = 0.9370000000000001
$ws.Range("F2").Value = 0.1095
$ws.Range("G2").Value = 0.0804

$ws.Range("B3").Value = 4.942
$ws.Range("C3").Value = 44.3776
$ws.Range("D3").Value = 6.6617
$ws.Range("E3").Value = 0.9473
$ws.Range("F3").Value = 0.1034
$ws.Range("G3").Value = 0.0819

$ws.Range("B4").Value = 3.9228
$ws.Range("C4").Value = 30.215
$ws.Range("D4").Value = 5.4968
$ws.Range("E4").Value = 0.9482
$ws.Range("F4").Value = 0.1096
$ws.Range("G4").Value = 0.078

$ws.Range("B5").Value = 5.4248
$ws.Range("C5").Value = 54.005
$ws.Range("D5").Value = 7.3488
$ws.Range("E5").Value = 0.9360000000000001
$ws.Range("F5").Value = 0.122
$ws.Range("G5").Value = 0.0895

$ws.Range("B6").Value = 7.3802
$ws.Range("C6").Value = 113.7728
$ws.Range("D6").Value = 10.6664
$ws.Range("E6").Value = 0.8529
$ws.Range("F6").Value = 0.1996
$ws.Range("G6").Value = 0.1543

$ws.Range("B7").Value = 5.3085
$ws.Range("C7").Value = 57.5679
$ws.Range("D7").Value = 7.3834
$ws.Range("E7").Value = 0.9243
$ws.Range("F7").Value = 0.1288
$ws.Range("G7").Value = 0.0968

$ws.Range("B8").Value = 1.1446
$ws.Range("C8").Value = 29.1212
$ws.Range("D8").Value = 1.7475
$ws.Range("E8").Value = 0.036
$ws.Range("F8").Value = 0.0359
$ws.Range("G8").Value = 0.029

# --- Sheet: pred_blend ---
$ws = $wb.Worksheets.Item("pred_blend")
$ws.Range("C2").Value = 4.9404
$ws.Range("D2").Value = 45.2628
$ws.Range("E2").Value = 6.7278
$ws.Range("F2").Value = 0.9052
$ws.Range("G2").Value = 0.107
$ws.Range("H2").Value = 0.0844

# --- Sheet: pred_stack ---
$ws = $wb.Worksheets.Item("pred_stack")
$ws.Range("C2").Value = 4.6081
$ws.Range("D2").Value = 38.7577
$ws.Range("E2").Value = 6.2256
$ws.Range("F2").Value = 0.9188
$ws.Range("G2").Value = 0.0997
$ws.Range("H2").Value = 0.07820000000000001

# --- Sheet: pred_final ---
$ws = $wb.Worksheets.Item("pred_final")
$ws.Range("B2").Value = "Stacking Regressor"
$ws.Range("C2").Value = 1.3449
$ws.Range("D2").Value = 4.1664
$ws.Range("E2").Value = 2.0412
$ws.Range("F2").Value = 0.9938
$ws.Range("G2").Value = 0.0344
$ws.Range("H2").Value = 0.0224